$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

# --- Title ---
Replace-Text "Unraveling the Enigma of Dark Matter" "Unveiling the Secrets of Life: A Journey Through Biology and Medicine"

# --- Author name ---
Replace-Text " Alexia Mitchell" " Lydia Mitchell"

# --- Email ---
Replace-Text "alexiamitchell@astros.ac.uk" "lydiamitchellphd@gmail.com"

Write-Host "done part 1"

# --- Body paragraph (Dark matter -> Biology/Medicine) ---
Replace-Text "Dark matter, an elusive and enigmatic component of the cosmos, has captivated the minds of scientists for decades" "The realm of biology and medicine stands as an expansive and intricate subject, inviting exploration into the very essence of life"

Replace-Text " This mysterious substance, invisible to our eyes and instruments, exerts a gravitational pull, shaping the structure and evolution of galaxies and galaxy clusters" " Biology delves into the mysteries of living organisms, from the smallest microorganisms to the grandest of creatures, unveiling the intricacies of their composition, function, and interdependence"

Write-Host "done part 2"

# Merge-and-replace: " Its existence...galaxies" + "." + " Understanding...evolution" -> new single sentence
Replace-Text " Its existence is inferred from its gravitational effects on visible matter, such as stars and gas, and from its influence on the motion of galaxies. Understanding the nature and properties of dark matter is a profound challenge in modern physics, with implications for our comprehension of the universe's composition, formation, and evolution" " Its intertwined partner, medicine, emerges as a beacon of healing and prevention, harnessing biological knowledge to address diseases, promote well-being, and ultimately prolong human existence"

Write-Host "done part 3"

Replace-Text "Unveiling the secrets of dark matter requires delving into the deepest mysteries of the universe" "Biology extends its reach into the microscopic realm, uncovering the wonders of cells, the basic units of life"

Replace-Text " From the subatomic realm to the vast cosmic web, scientists are employing a diverse array of techniques to probe this enigmatic entity" " Within these tiny compartments, intricate chemical processes orchestrated by DNA, the blueprint of life, dictate the characteristics and behaviors of organisms"

Replace-Text " Underground laboratories, sensitive observatories, and powerful particle accelerators are deployed in the quest to unravel the puzzle of dark matter's identity and its role in the grand symphony of the cosmos" " Biology illuminates the intricate mechanisms underpinning inheritance and evolution, revealing the astonishing diversity of life on Earth and the indissoluble link between all living things"

Replace-Text "The exploration of dark matter is a thrilling scientific endeavor that pushes the boundaries of human knowledge" "Medicine, drawing upon biological insights, emerges as a guardian of human health"

Replace-Text " With each new discovery, we inch closer to understanding the nature of this elusive substance and its profound impact on the universe" " It illuminates the causes and mechanisms of diseases, unraveling their complexities to pave the way for targeted interventions and therapies"

Write-Host "done part 4"

# --- Final sentence of body paragraph gets replaced, then many new sentences/line breaks appended ---
$rng = $d.Content
$found = $rng.Find.Execute(" As we unravel the enigma of dark matter, we illuminate the hidden workings of the cosmos, expanding our comprehension of the vast tapestry of reality that surrounds us", $true, $false, $false, $false, $false, $true, 1, $false, " The dedication of medical researchers and practitioners translates biological discoveries into life-saving treatments, vaccines, and preventive strategies", 2)
if (-not $found) { Write-Host "NOT FOUND: As we unravel the enigma..." }
$rng.Collapse(0)
$p = $rng.End

function InsertAt($text) {
    $r = $d.Range($p, $p)
    $before = $d.Content.End
    $r.InsertAfter($text)
    $after = $d.Content.End
    $p = $p + ($after - $before)
}

function BreakAt() {
    $r = $d.Range($p, $p)
    $before = $d.Content.End
    $r.InsertBreak(6)
    $after = $d.Content.End
    $p = $p + ($after - $before)
}

InsertAt "."
InsertAt " Moreover, medicine recognizes the influence of factors beyond biology, acknowledging the profound impact of social and environmental factors on health and well-being."
BreakAt
BreakAt
InsertAt "This intertwined journey through biology and medicine offers a lens to comprehend the intricacies of life, appreciate our place within the vast tapestry of existence, and harness our knowledge to improve human health. As we delve into this odyssey, may we cultivate a profound appreciation for the wonders of life and a commitment to preserving its delicate balance"

Write-Host "done part 5"
